# Apply the "math_L-curve" perturbation to the optimization_parameters sheet:
#  - drop the redundant extra "value" header cells (C1:F1)
#  - rename the "Model" parameter row to "production_function"
#  - add a new "L_curve" parameter row right after it (value = 1)
#  - drop the obsolete "Deletion" row from the Strain table
#  - leave the optimization_parameters sheet active, with the last row selected

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 1 only needs the A1/B1 header pair now; remove the duplicated "value" cells.
$ws.Range("C1:F1").ClearContents()

# "Model" -> "production_function" (row 8, same Sigmoid value stays in B8).
$ws.Range("A8").Value = "production_function"

# Insert the new "L_curve" row directly under it.
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 1
$ws.Range("B9").NumberFormat = "0.00E+00"

# The "Strain" table's "Deletion" row (now pushed down to row 17) is removed entirely.
$ws.Rows.Item(17).Delete()

# Make this the active sheet/tab, with the whole last row selected.
$ws.Activate()
$ws.Range("A17:XFD17").Select()
